# Ungroup "Group 5" on slide 1 (Diamond 3 + TextBox 4 annotation icon).
# PowerPoint's Ungroup() removes the wrapping <p:grpSp> and re-expresses
# each child shape's <a:xfrm> in slide coordinates, which is exactly the
# change recorded in the commit ("worked on annotation icons").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$grp = $s.Shapes.Item("Group 5")
$grp.Ungroup() | Out-Null
